# Generate Report for Handback
# Refresh the handback timestamps / status produced by the latest report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on the Overview sheet, and the matching
# "Correspond Handoff Datetime" column on the de-de sheet (same shared text).
$wsOverview.Range("G2").Value = "2016-11-03 20:26:51"
$wsOverview.Range("G3").Value = "2016-11-03 20:26:51"
$wsDeDe.Range("H2").Value = "2016-11-03 20:26:51"
$wsDeDe.Range("H3").Value = "2016-11-03 20:26:51"

# Status column ("ht" -> "mt") on both locale sheets.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# zh-cn "Correspond Handoff Datetime" / "Correspond Handback DateTime".
$wsZhCn.Range("H2").Value = "2016-11-03 20:26:37"
$wsZhCn.Range("H3").Value = "2016-11-03 20:26:37"
$wsZhCn.Range("K2").Value = "2016-11-03 20:27:28"
$wsZhCn.Range("K3").Value = "2016-11-03 20:27:28"

# de-de "Correspond Handback DateTime".
$wsDeDe.Range("K2").Value = "2016-11-03 20:27:47"
$wsDeDe.Range("K3").Value = "2016-11-03 20:27:47"
